$wb = $excel.ActiveWorkbook

# Sheet "展览" - row 8 (CM04 event) and row 10 (云蒸动漫音乐嘉年华)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 4621
$ws1.Range("F10").Value = 5054

# Sheet "全部类型" - same events, shifted one row down (row 9 and row 11)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 4621
$ws4.Range("F11").Value = 5054
